$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells we touch keep their original text (inline string) representation
# by forcing Text number format before assigning string values that Excel would
# otherwise auto-convert to numbers (e.g. "598.92", "1.00", "0.0000252").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.693.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.124.46'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.04%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.96'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.85%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.117.60'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.88%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.34'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.465'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.95'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.637.45'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.90%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.778.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.120.09'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.80'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.25'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.67'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.705'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.58'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.65%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.33'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.60%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.22'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.94'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.77%  '

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.97'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.50%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.55%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.41%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.56'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0737'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.92'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -7.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '436.34'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.66%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.57%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.24'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.863.05'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.258'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.19'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.07%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.78'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.79%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.02%  '
